# Build site at 2023-04-12 14:53:07 UTC
# Applies the content update to disciplina 8800010 (Canto Coral II) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 10 (Objetivos / B+C): replace teacher-name text with the real
#    Portuguese objectives paragraph.
# ---------------------------------------------------------------------------
$objetivosPt = '1. Aproximar o aluno do seu aparelho vocal, ao nível da expressão falada ou cantada, em sua expressão individual ou coletiva (coral). 2. Propiciar ao aluno o contato com o repertório coral e seu papel no desenvolvimento da linguagem musical. 3. Propiciar ao aluno a oportunidade de fazer música vocal em conjunto, com o aprendizado técnico de parâmetros como: afinação, precisão, equilíbrio, fraseado etc.'
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# ---------------------------------------------------------------------------
# 2. Insert a new row 13 (the "Docentes responsáveis" value row) before the
#    current "Programa resumido" row, copying formatting from the row
#    below so row heights of the following rows are preserved unchanged.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert([System.Reflection.Missing]::Value, -4163)
$ws.Rows.Item(13).Clear()

$docente = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# Fix up the new cells' styles (column B/C body styles) by pasting the
# formats from the untouched reference row 11.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row 14 (was 13, "Programa resumido:"): replace "Semestral" with the
#    short Portuguese syllabus paragraph.
# ---------------------------------------------------------------------------
$programaResumidoPt = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto. Colocação da emissão em "Bocca Chiusa". O canto coral em uníssono. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral.'
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt

# ---------------------------------------------------------------------------
# 4. Row 16 (was 15, "Programa:"): replace stray "01/01/2017" with the full
#    Portuguese programa paragraph.
# ---------------------------------------------------------------------------
$programaPt = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto - exercícios para localização da respiração baixa e média. Apoio e coluna de ar. Colocação da emissão em "Bocca Chiusa" relaxamento da mandíbula, posicionamento da língua, suspensão do palato mole, exercícios de percepção do local onde a voz está se colocando, conexão do apoio e emissão vocal, passagem da "Bocca Chiusa" para vogais e outros sons nasais e guturais. O canto coral em uníssono: afinação, uniformidade tímbrica, precisão rítmica. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral: testagem da escuta harmônica e afinação à 1ª vista a várias vozes, memorização, fundamentos de teoria musical. Montagem e aperfeiçoamento de peças musicais - promovendo a aplicação das técnicas aprendidas. Conexão entre diafragma e emissão vocal.'
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# ---------------------------------------------------------------------------
# 5. Row 19 (was 18, "Método:"): replace teacher-name text with the
#    "A cada semestre..." paragraph.
# ---------------------------------------------------------------------------
$metodoTxt = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Range("B19").Value = $metodoTxt
$ws.Range("C19").Value = $metodoTxt

# ---------------------------------------------------------------------------
# 6. Row 20 (was 19, "Critério:"): replace "A cada semestre..." with the
#    "Sendo uma atividade prática..." paragraph.
# ---------------------------------------------------------------------------
$criterioTxt = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Range("B20").Value = $criterioTxt
$ws.Range("C20").Value = $criterioTxt

# ---------------------------------------------------------------------------
# 7. Row 21 (was 20, "Norma de recuperação:"): replace the "Sendo uma
#    atividade prática..." text with "Não tem".
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "Não tem"
$ws.Range("C21").Value = "Não tem"

# ---------------------------------------------------------------------------
# 8. Row 22 (was 21, "Bibliografia:"): replace "Não tem" with the full
#    bibliography text.
# ---------------------------------------------------------------------------
$bibliografiaTxt = 'BEHLAU, Mara; REHDER, Maria Inês. Higiene vocal para o Canto Coral. Rio de Janeiro: Revinter, 1997.COELHO, Helena Wöhl. Técnica vocal para coros. 7ª Edição. São Leopoldo: Sinodal, 1994.CORBIN, Lynn Ann. Vocal pedagogy in the choral rehearsal: The selected concepts on choral tone quality, understanding of the singing process, and attitudes toward choir participation. 1982. 119 f. Tese (Doutorado em Filosofia) - Graduate School of The Ohio State University, Ohio State University, 1982.COSTA, Paulo Rubens Moraes. Diagnose em Canto Coral: parâmetros para análise e ferramentas para a avaliação. 2005. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2005.DRAHAN, Snizhana. Ouvir a voz: a percepção da produção vocal pelo regente coral - método eformação. 2007. 146 f. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2007.FERNANDES, Angelo José. O regente moderno e a construção da sonoridade coral: uma metodologia de preparo vocal para coros. 2009. 475 f. Tese (Doutorado em Música) - Instituto de Artes, Universidade Federal de Campinas, Campinas, 2009.HERR, Martha. Considerações para a classificação da voz do coralista. In: FERREIRA, LésliePiccolotto et al. Voz profissional: o profissional da voz. Carapicuíba: Pró-fono DepartamentoEditorial, 1995. p. 51-56.MILLER, Richard. The structure of singing. Boston, Massachusetts: Schirmer, 1986.RAMOS, Marco Antonio da Silva. Canto Coral: do repertório temático à construção do programa. 1988. 492f. Dissertação (Mestrado em Artes) - Departamento de Biblioteconomia e Documentação, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 1988.______. Memorial. 2011. Memorial (Professor titular) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2011.______. O ensino da regência coral. 2003. 107f. Tese (Livre Docência) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2003.VIDEIRA JR., Mário Rodrigues. Educação musical através do coro: A experiência no Projeto Comunicantus da ECA - USP. 2001. 65 f. Trabalho de Conclusão de Curso (Licenciatura em Música) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2001.'
$ws.Range("B22").Value = $bibliografiaTxt
$ws.Range("C22").Value = $bibliografiaTxt

Write-Host "Edit complete"
